# Vocab augmented (author = Alex)
#
# Adds a new column J ("GloVe used") to the existing comparison table and
# appends a long sweep of nb_iters values (50..50000) for both "their" and
# "ours" GloVe embeddings (rows 16-33), replacing the previous two
# trailing rows (16-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    # Force the cell to be stored as text (matches source cells that were
    # already shared-string/text, e.g. "0.996" and "68.5"), then drop the
    # temporary number-format override so no stray style sticks around.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- New header cell (bold, matching the rest of the header row) -------
$ws.Cells.Item(5, 10).Value = "GloVe used"
$ws.Cells.Item(5, 10).Font.Bold = $true

# --- Rows 6-15: existing data rows simply gain a "-" (not applicable) in
#     the new "GloVe used" column ---------------------------------------
for ($r = 6; $r -le 15; $r++) {
    $ws.Cells.Item($r, 10).Value = "-"
}

# --- Rows 16-33: replace the old two rows with a full nb_iters sweep for
#     "their" GloVe (rows 16-24) and "ours" GloVe (rows 25-33) ----------
$nbIters = @(50, 100, 500, 1000, 2000, 5000, 10000, 20000, 50000)

$row = 16
foreach ($n in $nbIters) {
    $ws.Cells.Item($row, 3).Value = "neural network"
    $ws.Cells.Item($row, 4).Value = 35
    $ws.Cells.Item($row, 5).Value = "default"
    $ws.Cells.Item($row, 6).Value = "default"
    Set-TextValue $row 7 "0.996"
    $ws.Cells.Item($row, 8).Value = $n
    Set-TextValue $row 9 "68.5"
    $ws.Cells.Item($row, 10).Value = "their"
    $row++
}

foreach ($n in $nbIters) {
    $ws.Cells.Item($row, 3).Value = "neural network"
    $ws.Cells.Item($row, 4).Value = 35
    $ws.Cells.Item($row, 5).Value = "default"
    $ws.Cells.Item($row, 6).Value = "default"
    Set-TextValue $row 7 "0.996"
    $ws.Cells.Item($row, 8).Value = $n
    Set-TextValue $row 9 "68.5"
    $ws.Cells.Item($row, 10).Value = "ours"
    $row++
}

# --- Refresh the view: scroll so row 33 is the top-left visible row and
#     select the cell just below the new data (mirrors the author's
#     on-save selection state) -------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Range("I52").Select()
